# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two "section header" rows (row 5 "situação do domicílio" and
# row 8 "grandes regiões e unidades da federação") that carried a label but
# no data - a leftover from the pandas export. Removing those rows shifts
# every row below them up, realigning each label with the data row that
# actually belongs to it. The column-2 header label is also fixed from the
# stray pandas placeholder "unnamed: 1_level_1" to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-labelled sub-header in row 2 (B2).
$ws.Cells.Item(2, 2).Value = "total"

# Remove the "situação do domicílio" section-header row (row 5). Excel
# shifts rows 6+ up by one, so what was row 6 ("urbana") now lands on row 5,
# etc.
$ws.Rows(5).Delete()

# After the first delete, the "grandes regiões e unidades da federação"
# section-header row (originally row 8) is now row 7. Remove it too, which
# shifts the region/state rows up by one more.
$ws.Rows(7).Delete()
